$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B header: "Age" -> "Birthday"
$ws.Range("B1").Value = "Birthday"

# Column B values: ages -> birthdate 1/1/1980 (serial 29221), formatted as a short date
$ws.Range("B2").Value = 29221
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy($ws.Range("B3"))
$ws.Range("B2").Copy($ws.Range("B4"))
$ws.Range("B2").Copy($ws.Range("B5"))

# Update the selected cell to match the author's final cursor position
$ws.Range("E9").Select()
